$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Clear()
$ws.Range("C13:C30").Clear()

$ws.Range("B11").Value = "Wed Mar 20 23:02:17 EDT 2024"
$ws.Range("B12").Value = "Wed Mar 20 23:02:27 EDT 2024"

$ws.Range("C13:C30").Select() | Out-Null
